$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap out the album-cover image paths referenced by the "Gruppe 16" and
# "Gruppe 18" rows (the shared-string table entries for B15 / B17).
$ws.Range("B15").Value = "/album_covers/artist_16.jpg"
$ws.Range("B17").Value = "/album_covers/gruppe_18.png"

# Move the selection to B18 and scroll the window so row 4 is the top
# visible row (sheetView topLeftCell="A4").
$ws.Range("B18").Select()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
